$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1935483870967742
$ws.Range("C2").Value = 0.5833333333333334
$ws.Range("J2").Value = 0.01881720430107527
$ws.Range("P2").Value = 0.1344086021505376
$ws.Range("S2").Value = 0.06989247311827956
# Row 3
$ws.Range("B3").Value = 0.008733624454148471
$ws.Range("C3").Value = 0.03930131004366812
$ws.Range("J3").Value = 0.008733624454148471
$ws.Range("P3").Value = 0.7685589519650655
$ws.Range("S3").Value = 0.1746724890829694
# Row 4
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.7083333333333334
$ws.Range("S4").Value = 0.2083333333333333
# Row 6
$ws.Range("B6").Value = 0.05603448275862069
$ws.Range("D6").Value = 0.02155172413793104
$ws.Range("F6").Value = 0.05603448275862069
$ws.Range("J6").Value = 0.2801724137931034
$ws.Range("O6").Value = 0.02586206896551724
$ws.Range("Q6").Value = 0.1767241379310345
$ws.Range("R6").Value = 0.04310344827586207
$ws.Range("S6").Value = 0.3405172413793103
# Row 7
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.01075268817204301
$ws.Range("F7").Value = 0.04659498207885305
$ws.Range("J7").Value = 0.1362007168458781
$ws.Range("O7").Value = 0.01433691756272401
$ws.Range("Q7").Value = 0.1792114695340502
$ws.Range("R7").Value = 0.05734767025089606
$ws.Range("S7").Value = 0.4444444444444444
# Row 8
$ws.Range("B8").Value = 0.09569377990430622
$ws.Range("D8").Value = 0.02392344497607655
$ws.Range("E8").Value = 0.002392344497607655
$ws.Range("F8").Value = 0.0645933014354067
$ws.Range("J8").Value = 0.1220095693779904
$ws.Range("O8").Value = 0.03349282296650718
$ws.Range("Q8").Value = 0.1746411483253588
$ws.Range("R8").Value = 0.0861244019138756
$ws.Range("S8").Value = 0.3971291866028708
# Row 9
$ws.Range("B9").Value = 0.1185567010309278
$ws.Range("D9").Value = 0.005154639175257732
$ws.Range("F9").Value = 0.04639175257731959
$ws.Range("J9").Value = 0.1134020618556701
$ws.Range("O9").Value = 0.005154639175257732
$ws.Range("Q9").Value = 0.1907216494845361
$ws.Range("R9").Value = 0.07216494845360824
$ws.Range("S9").Value = 0.4484536082474227
# Row 10
$ws.Range("B10").Value = 0.1277346506704305
$ws.Range("D10").Value = 0.0218772053634439
$ws.Range("F10").Value = 0.05363443895553988
$ws.Range("J10").Value = 0.1171489061397318
$ws.Range("O10").Value = 0.01976005645730416
$ws.Range("Q10").Value = 0.2180663373323924
$ws.Range("R10").Value = 0.06422018348623854
$ws.Range("S10").Value = 0.3775582215949189
# Row 11
$ws.Range("G11").Value = 0.1570438799076213
$ws.Range("J11").Value = 0.09699769053117784
$ws.Range("K11").Value = 0.1939953810623557
$ws.Range("L11").Value = 0.5450346420323325
$ws.Range("S11").Value = 0.006928406466512702
# Row 12
$ws.Range("G12").Value = 0.7302904564315352
$ws.Range("J12").Value = 0.1991701244813278
$ws.Range("K12").Value = 0.008298755186721992
$ws.Range("L12").Value = 0.02489626556016597
$ws.Range("S12").Value = 0.03734439834024896
# Row 13
$ws.Range("G13").Value = 0.7666666666666667
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.03333333333333333
# Row 15
$ws.Range("F15").Value = 0.02845528455284553
$ws.Range("H15").Value = 0.1260162601626016
$ws.Range("I15").Value = 0.07317073170731707
$ws.Range("J15").Value = 0.3455284552845528
$ws.Range("K15").Value = 0.04471544715447155
$ws.Range("M15").Value = 0.02439024390243903
$ws.Range("O15").Value = 0.1016260162601626
$ws.Range("S15").Value = 0.2560975609756098
# Row 16
$ws.Range("F16").Value = 0.03174603174603174
$ws.Range("H16").Value = 0.1388888888888889
$ws.Range("I16").Value = 0.07539682539682539
$ws.Range("J16").Value = 0.4246031746031746
$ws.Range("K16").Value = 0.1468253968253968
$ws.Range("M16").Value = 0.0119047619047619
$ws.Range("O16").Value = 0.04365079365079365
$ws.Range("S16").Value = 0.126984126984127
# Row 17
$ws.Range("F17").Value = 0.0398406374501992
$ws.Range("H17").Value = 0.1693227091633466
$ws.Range("I17").Value = 0.07370517928286853
$ws.Range("J17").Value = 0.398406374501992
$ws.Range("K17").Value = 0.1254980079681275
$ws.Range("M17").Value = 0.02191235059760956
$ws.Range("N17").Value = 0.00199203187250996
$ws.Range("O17").Value = 0.05776892430278884
$ws.Range("S17").Value = 0.1115537848605578
# Row 18
$ws.Range("F18").Value = 0.02409638554216868
$ws.Range("H18").Value = 0.1325301204819277
$ws.Range("I18").Value = 0.06626506024096386
$ws.Range("J18").Value = 0.4819277108433735
$ws.Range("K18").Value = 0.1566265060240964
$ws.Range("M18").Value = 0.03012048192771084
$ws.Range("O18").Value = 0.04819277108433735
$ws.Range("S18").Value = 0.06024096385542169
# Row 19
$ws.Range("F19").Value = 0.02551381998582565
$ws.Range("H19").Value = 0.1750531537916371
$ws.Range("I19").Value = 0.07795889440113395
$ws.Range("J19").Value = 0.3586109142452162
$ws.Range("K19").Value = 0.1438695960311835
$ws.Range("M19").Value = 0.02693125442948264
$ws.Range("N19").Value = 0.0007087172218284905
$ws.Range("O19").Value = 0.06024096385542169
$ws.Range("S19").Value = 0.1311126860382707
